$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) - mirror the existing header formatting (bold/border/center/top)
# by copying the format from an existing header cell (H1) rather than rebuilding it by hand,
# so the same underlying cell style is reused.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J for rows 2-14
$data = @{
    2  = @(6, 6)
    3  = @(7, 8)
    4  = @(6, 6)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(7, 7)
    8  = @(7, 7)
    9  = @(5, 6)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(1, 3)
    14 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
